$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumo_por_Cliente")

# Atualizacao dos dados da BIBI: coluna "situacao" (J) recalculada para 18 clientes
# (meses sem comprar avancaram conforme a nova data de referencia do calculo)
$ws.Range("J4").Value = "INATIVO - 35.0 meses sem comprar"
$ws.Range("J9").Value = "INATIVO - 18.0 meses sem comprar"
$ws.Range("J10").Value = "INATIVO - 1.5 meses sem comprar"
$ws.Range("J11").Value = "INATIVO - 3.9 meses sem comprar"
$ws.Range("J18").Value = "INATIVO - 10.3 meses sem comprar"
$ws.Range("J20").Value = "INATIVO - 36.2 meses sem comprar"
$ws.Range("J33").Value = "INATIVO - 13.5 meses sem comprar"
$ws.Range("J34").Value = "INATIVO - 25.7 meses sem comprar"
$ws.Range("J38").Value = "INATIVO - 33.3 meses sem comprar"
$ws.Range("J45").Value = "INATIVO - 1.4 meses sem comprar"
$ws.Range("J49").Value = "INATIVO - 8.9 meses sem comprar"
$ws.Range("J64").Value = "INATIVO - 20.8 meses sem comprar"
$ws.Range("J66").Value = "INATIVO - 11.5 meses sem comprar"
$ws.Range("J72").Value = "INATIVO - 20.8 meses sem comprar"
$ws.Range("J84").Value = "INATIVO - 8.4 meses sem comprar"
$ws.Range("J89").Value = "INATIVO - 14.0 meses sem comprar"
$ws.Range("J94").Value = "INATIVO - 18.1 meses sem comprar"
$ws.Range("J98").Value = "INATIVO - 21.8 meses sem comprar"
